$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated (Latitude, Longitude) readings for rows 2-9 (columns B and C)
$ws.Range("B2").Value = 42.989040840000001
$ws.Range("C2").Value = -81.228214929999993

$ws.Range("B3").Value = 42.986799697825198
$ws.Range("C3").Value = -81.227937513089998

$ws.Range("B4").Value = 42.986657532755402
$ws.Range("C4").Value = -81.230996674159698

$ws.Range("B5").Value = 42.984427370468701
$ws.Range("C5").Value = -81.230586739052399

$ws.Range("B6").Value = 42.984427370468701
$ws.Range("C6").Value = -81.230586739052399

$ws.Range("B7").Value = 42.984427370468701
$ws.Range("C7").Value = -81.230586739052399

$ws.Range("B8").Value = 42.984427370468701
$ws.Range("C8").Value = -81.230586739052399

$ws.Range("B9").Value = 42.984395306414903
$ws.Range("C9").Value = -81.233651596635298

# Selection moved from U4 to A5
$ws.Range("A5").Select()
